$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C2:C6) from 2023-09-01 (45170)
# to 2023-09-05 (45174), keeping the existing date formatting/style.
$ws.Range("C2:C6").Value = 45174
